$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("F74").Value = 38
$ws.Range("G74").Value = 4374.56
$ws.Range("F77").Value = 268
$ws.Range("G77").Value = 17071.6
$ws.Range("F79").Value = 63
$ws.Range("G79").Value = 4975.74
$ws.Range("F92").Value = 38
$ws.Range("G92").Value = 4767.86
$ws.Range("B96").Value = 153056.43
$ws.Range("F106").Value = 5
$ws.Range("G106").Value = 18193.5
$ws.Range("F108").Value = 4
$ws.Range("G108").Value = 11283.64
$ws.Range("B114").Value = 53818.15
$ws.Range("F118").Value = 211
$ws.Range("G118").Value = 10984.66
$ws.Range("F122").Value = 75
$ws.Range("G122").Value = 4925.25
$ws.Range("B124").Value = 53755.2
$ws.Range("F144").Value = 23
$ws.Range("G144").Value = 11905.26
$ws.Range("F151").Value = 19
$ws.Range("G151").Value = 6959.51
$ws.Range("F152").Value = 100
$ws.Range("G152").Value = 10881
$ws.Range("F153").Value = 7
$ws.Range("G153").Value = 12815.46
$ws.Range("F154").Value = 4
$ws.Range("G154").Value = 10598.92
$ws.Range("F156").Value = 5
$ws.Range("G156").Value = 27967.75
$ws.Range("B159").Value = 222340.93
$ws.Range("F184").Value = 310
$ws.Range("G184").Value = 2476.9
$ws.Range("B186").Value = 10281.6
$ws.Range("F188").Value = 181
$ws.Range("G188").Value = 11728.8
$ws.Range("F190").Value = 79
$ws.Range("G190").Value = 6863.52
$ws.Range("F191").Value = 52
$ws.Range("G191").Value = 4591.08
$ws.Range("B195").Value = 24695.19
$ws.Range("F229").Value = 8
$ws.Range("G229").Value = 2430
$ws.Range("B230").Value = 24382.12
$ws.Range("F269").Value = 30
$ws.Range("G269").Value = 2631
$ws.Range("B271").Value = 25579.84
$ws.Range("B284").Value = 57802
$ws.Range("E284").Value = 162.71
$ws.Range("F284").Value = -79
$ws.Range("G284").Value = -11334.92
$ws.Range("B285").Value = 63531
$ws.Range("E285").Value = 152.53
$ws.Range("F285").Value = 45
$ws.Range("G285").Value = 6456.6
$ws.Range("F289").Value = 28
$ws.Range("G289").Value = 1436.96
$ws.Range("F310").Value = 478
$ws.Range("G310").Value = 81895.74000000001
$ws.Range("B315").Value = 151211.01
$ws.Range("F353").Value = 7
$ws.Range("G353").Value = 499.52
$ws.Range("B363").Value = 147461.66
$ws.Range("F389").Value = 21
$ws.Range("G389").Value = 1441.65
$ws.Range("F396").Value = 10
$ws.Range("G396").Value = 2515
$ws.Range("B399").Value = 22971.56
$ws.Range("F403").Value = 181
$ws.Range("G403").Value = 9379.42
$ws.Range("F414").Value = 13
$ws.Range("G414").Value = 798.33
$ws.Range("B415").Value = 21056.09
$ws.Range("F430").Value = 7
$ws.Range("G430").Value = 23252.25
$ws.Range("B443").Value = 200476.06
$ws.Range("B462").Value = 45709
$ws.Range("E462").Value = 15.69
$ws.Range("F462").Value = -300
$ws.Range("G462").Value = -3945
$ws.Range("B463").Value = 64925
$ws.Range("E463").Value = 13.97
$ws.Range("F463").Value = 111
$ws.Range("G463").Value = 1459.65
$ws.Range("B540").Value = 64810
$ws.Range("E540").Value = 291.22
$ws.Range("F540").Value = 4
$ws.Range("G540").Value = 1095.68
$ws.Range("B541").Value = 53319
$ws.Range("E541").Value = 310.64
$ws.Range("F541").Value = -6
$ws.Range("G541").Value = -1643.52
$ws.Range("B571").Value = 60022
$ws.Range("E571").Value = 37.22
$ws.Range("F571").Value = -113
$ws.Range("G571").Value = -3709.79
$ws.Range("B572").Value = 64830
$ws.Range("E572").Value = 34.9
$ws.Range("F572").Value = 101
$ws.Range("G572").Value = 3315.83
$ws.Range("F616").Value = 14
$ws.Range("G616").Value = 23671.34
$ws.Range("F617").Value = 5
$ws.Range("G617").Value = 11380.4
$ws.Range("B621").Value = 171855.5
$ws.Range("F624").Value = 511
$ws.Range("G624").Value = 62265.35
$ws.Range("F626").Value = 79
$ws.Range("G626").Value = 6759.24
$ws.Range("F628").Value = 73
$ws.Range("G628").Value = 6050.97
$ws.Range("B631").Value = 113040.19
$ws.Range("F665").Value = 1
$ws.Range("G665").Value = 4269.07
$ws.Range("F666").Value = 6
$ws.Range("G666").Value = 28686.36
$ws.Range("F667").Value = 2
$ws.Range("G667").Value = 10586.68
$ws.Range("F672").Value = 7
$ws.Range("G672").Value = 19126.66
$ws.Range("F673").Value = 56
$ws.Range("G673").Value = 13746.32
$ws.Range("B676").Value = 216536.3
$ws.Range("F694").Value = 176
$ws.Range("G694").Value = 14354.56
$ws.Range("F699").Value = 477
$ws.Range("G699").Value = 63488.7
$ws.Range("F705").Value = 104
$ws.Range("G705").Value = 14041.04
$ws.Range("F706").Value = 716
$ws.Range("G706").Value = 86428.36
$ws.Range("B707").Value = 179269.86
$ws.Range("B719").Value = 65362
$ws.Range("F719").Value = 18
$ws.Range("G719").Value = 735.66
$ws.Range("B720").Value = 65079
$ws.Range("F720").Value = 6
$ws.Range("G720").Value = 245.22
$ws.Range("F727").Value = 34
$ws.Range("G727").Value = 1937.66
$ws.Range("B730").Value = 9118.450000000001
$ws.Range("F749").Value = 3420
$ws.Range("G749").Value = 557836.2
$ws.Range("F751").Value = 282
$ws.Range("G751").Value = 79769.34
$ws.Range("F752").Value = 250
$ws.Range("G752").Value = 36162.5
$ws.Range("B756").Value = 674616.09
$ws.Range("B775").Value = 5200212.15
$ws.Range("B776").Value = 5200212.15

Write-Output "Applied 153 cell updates"